$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price (D) column stores numeric-looking strings (e.g. "43.46")
# as plain text in the source workbook (t="inlineStr", no numFmt).
# Excel's COM Value setter auto-converts such strings to numbers, so
# we briefly force Text number-format before assigning, then restore
# the cell's style to Normal (no lingering numFmt) afterwards so only
# the text content differs from the original, not the cell style.
$priceCells = @("D2", "D3", "D5", "D6", "D7", "D9", "D11", "D12", "D13", "D14", "D15", "D16", "D18", "D20", "D21", "D23", "D24", "D26", "D27", "D28", "D29", "D30", "D34", "D36", "D38", "D39", "D42", "D44", "D46", "D47", "D48", "D49")
foreach ($cellRef in $priceCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

# Apply the updated coin data (price, volume-change, and for two rows
# that swapped ranking position: coin name + link as well).
$ws.Range("D2").Value = '67.337.22'
$ws.Range("E2").Value = '  +1.02%  '
$ws.Range("D3").Value = '3.947.93'
$ws.Range("E3").Value = '  +4.09%  '
$ws.Range("E4").Value = '  -0.01%  '
$ws.Range("D5").Value = '471.45'
$ws.Range("E5").Value = '  +8.92%  '
$ws.Range("D6").Value = '145.80'
$ws.Range("E6").Value = '  +4.06%  '
$ws.Range("D7").Value = '0.624'
$ws.Range("E7").Value = '  -0.06%  '
$ws.Range("D9").Value = '0.733'
$ws.Range("E9").Value = '  -0.23%  '
$ws.Range("E10").Value = '  +7.72%  '
$ws.Range("D11").Value = '0.0000338'
$ws.Range("E11").Value = '  +6.91%  '
$ws.Range("D12").Value = '43.46'
$ws.Range("D13").Value = '4.571.58'
$ws.Range("E13").Value = '  +3.62%  '
$ws.Range("D14").Value = '10.39'
$ws.Range("E14").Value = '  -0.73%  '
$ws.Range("D15").Value = '15.19'
$ws.Range("E15").Value = '  +1.06%  '
$ws.Range("D16").Value = '3.889.14'
$ws.Range("E16").Value = '  +1.37%  '
$ws.Range("E17").Value = '  -0.27%  '
$ws.Range("D18").Value = '19.85'
$ws.Range("E18").Value = '  -0.39%  '
$ws.Range("E19").Value = '  +2.14%  '
$ws.Range("D20").Value = '67.584.07'
$ws.Range("E20").Value = '  +1.22%  '
$ws.Range("D21").Value = '433.48'
$ws.Range("E21").Value = '  +5.68%  '
$ws.Range("E22").Value = '  +3.18%  '
$ws.Range("D23").Value = '14.54'
$ws.Range("E23").Value = '  -0.62%  '
$ws.Range("D24").Value = '87.39'
$ws.Range("E24").Value = '  +2.56%  '
$ws.Range("E25").Value = '  +7.79%  '
$ws.Range("D26").Value = '38.63'
$ws.Range("E26").Value = '  +5.09%  '
$ws.Range("B27").Value = 'Filecoin'
$ws.Range("C27").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D27").Value = '10.30'
$ws.Range("E27").Value = '  +5.18%  '
$ws.Range("B28").Value = 'LEO'
$ws.Range("C28").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D28").Value = '5.75'
$ws.Range("E28").Value = '  +2.53%  '
$ws.Range("D29").Value = '9.63'
$ws.Range("E29").Value = '  -0.82%  '
$ws.Range("D30").Value = '726.53'
$ws.Range("E30").Value = '  +0.50%  '
$ws.Range("E31").Value = '  -2.22%  '
$ws.Range("E32").Value = '  -2.37%  '
$ws.Range("E33").Value = '  +4.00%  '
$ws.Range("D34").Value = '43.15'
$ws.Range("E34").Value = '  +4.34%  '
$ws.Range("E35").Value = '  +1.83%  '
$ws.Range("D36").Value = '57.93'
$ws.Range("E36").Value = '  +3.71%  '
$ws.Range("E37").Value = '  -0.03%  '
$ws.Range("D38").Value = '0.0₃0780'
$ws.Range("E38").Value = '  +13.40%  '
$ws.Range("D39").Value = '5.38'
$ws.Range("E39").Value = '  -5.43%  '
$ws.Range("E40").Value = '  +0.73%  '
$ws.Range("E41").Value = '  +1.54%  '
$ws.Range("D42").Value = '2.59'
$ws.Range("E42").Value = '  -5.83%  '
$ws.Range("E43").Value = '  -0.30%  '
$ws.Range("D44").Value = '0.337'
$ws.Range("E44").Value = '  +4.52%  '
$ws.Range("E45").Value = '  -0.21%  '
$ws.Range("D46").Value = '2.82'
$ws.Range("E46").Value = '  +4.67%  '
$ws.Range("B47").Value = 'LidoDAOToken'
$ws.Range("C47").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D47").Value = '3.45'
$ws.Range("E47").Value = '  +3.89%  '
$ws.Range("B48").Value = 'ARBITRUM'
$ws.Range("C48").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D48").Value = '2.19'
$ws.Range("E48").Value = '  +4.75%  '
$ws.Range("D49").Value = '147.41'
$ws.Range("E49").Value = '  +3.75%  '
$ws.Range("E50").Value = '  -2.05%  '
$ws.Range("E51").Value = '  +1.63%  '

foreach ($cellRef in $priceCells) {
    $ws.Range($cellRef).Style = "Normal"
}
